# Apply updated crypto price/volume figures scraped by GitHub Actions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay text even if it looks like a number
# (Excel would otherwise silently convert "580.05" etc. into a float).
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    if ($text -match "^[+-]?[0-9]*\.?[0-9]+$") {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

Set-TextValue 'D2' '61.885.60'
Set-TextValue 'E2' '  +4.53%  '
Set-TextValue 'D3' '3.080.22'
Set-TextValue 'E3' '  +3.19%  '
Set-TextValue 'E4' '  +0.06%  '
Set-TextValue 'D5' '580.05'
Set-TextValue 'E5' '  +3.13%  '
Set-TextValue 'D6' '142.32'
Set-TextValue 'E6' '  +2.79%  '
Set-TextValue 'E7' '  -0.03%  '
Set-TextValue 'D8' '3.069.93'
Set-TextValue 'E8' '  +3.33%  '
Set-TextValue 'E9' '  +1.18%  '
Set-TextValue 'E10' '  +5.64%  '
Set-TextValue 'D11' '5.74'
Set-TextValue 'E11' '  +11.34%  '
Set-TextValue 'E12' '  +2.83%  '
Set-TextValue 'E13' '  +4.82%  '
Set-TextValue 'D14' '35.35'
Set-TextValue 'E14' '  +5.00%  '
Set-TextValue 'E15' '  +0.30%  '
Set-TextValue 'D16' '3.590.61'
Set-TextValue 'E16' '  +3.25%  '
Set-TextValue 'D17' '7.26'
Set-TextValue 'E17' '  +1.26%  '
Set-TextValue 'D18' '3.079.47'
Set-TextValue 'E18' '  +3.29%  '
Set-TextValue 'D19' '61.820.77'
Set-TextValue 'E19' '  +4.57%  '
Set-TextValue 'D20' '447.95'
Set-TextValue 'E20' '  +4.67%  '
Set-TextValue 'E21' '  +2.42%  '
Set-TextValue 'E22' '  +2.35%  '
Set-TextValue 'D23' '7.44'
Set-TextValue 'E23' '  +4.86%  '
Set-TextValue 'D24' '13.81'
Set-TextValue 'E24' '  +3.34%  '
Set-TextValue 'D25' '81.91'
Set-TextValue 'E25' '  +1.38%  '
Set-TextValue 'E26' '  +0.13%  '
Set-TextValue 'E27' '  +4.99%  '
Set-TextValue 'E28' '  +0.15%  '
Set-TextValue 'B29' 'PancakeSwap'
Set-TextValue 'C29' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D29' '2.67'
Set-TextValue 'E29' '  +5.06%  '
Set-TextValue 'B30' 'RenderToken'
Set-TextValue 'C30' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D30' '8.26'
Set-TextValue 'E30' '  +7.31%  '
Set-TextValue 'D31' '6.82'
Set-TextValue 'E31' '  +12.19%  '
Set-TextValue 'E32' '  +12.95%  '
Set-TextValue 'D33' '26.83'
Set-TextValue 'E33' '  +4.49%  '
Set-TextValue 'D34' '1.04'
Set-TextValue 'E34' '  +4.60%  '
Set-TextValue 'E35' '  +3.06%  '
Set-TextValue 'E36' '  +3.41%  '
Set-TextValue 'E37' '  +5.65%  '
Set-TextValue 'D38' '50.12'
Set-TextValue 'E38' '  +2.05%  '
Set-TextValue 'D39' '2.98'
Set-TextValue 'E39' '  +9.85%  '
Set-TextValue 'D40' '8.79'
Set-TextValue 'E40' '  +1.79%  '
Set-TextValue 'D41' '421.41'
Set-TextValue 'E41' '  +5.19%  '
Set-TextValue 'E42' '  +5.68%  '
Set-TextValue 'D43' '2.910.25'
Set-TextValue 'E43' '  +5.17%  '
Set-TextValue 'E44' '  +9.50%  '
Set-TextValue 'E45' '  +1.18%  '
Set-TextValue 'E46' '  +7.15%  '
Set-TextValue 'B47' 'USDe'
Set-TextValue 'C47' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D47' '0.999'
Set-TextValue 'E47' '  +0.01%  '
Set-TextValue 'B48' 'Arweave'
Set-TextValue 'C48' 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue 'D48' '35.12'
Set-TextValue 'E48' '  +3.27%  '
Set-TextValue 'E49' '  +2.18%  '
Set-TextValue 'E50' '  +0.84%  '
Set-TextValue 'D51' '24.20'
Set-TextValue 'E51' '  +3.47%  '
